# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
#
# A new weekly record is inserted into the price table right after the
# existing row 447, pushing the former rows 448-473 down to 449-474
# (the sheet's used range grows from A1:R473 to A1:R474).
#
# The new row (448) contains:
#   Mercado ID=10, Mercado="Vega Modelo de Temuco", Region="La Araucania",
#   Fecha=2023-04-25 (serial 45041), Codreg=9, Categoria ID=100112044,
#   Categoria="Perejil", Variedad="Sin especificar", Calidad="Primera",
#   Volumen=30, Precio minimo=4000, Precio maximo=4000,
#   Precio promedio ponderado=4000, Unidad="$/docena de atados (3 kilos)",
#   Origen="Provincia de Cautin", Precio $/Kg=1333, Kg o Unidades=3,
#   Clasificacion="Hortaliza".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 448; this shifts rows 448:473 down to
# 449:474 and carries the formatting (incl. the date number format of
# column D) down from the row above, just like Excel's normal
# "Insert Sheet Rows" behaviour.
$ws.Rows.Item(448).Insert()

$newRow = 448

$ws.Cells.Item($newRow, 1).Value2 = 10
$ws.Cells.Item($newRow, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value2 = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value2 = 45041
$ws.Cells.Item($newRow, 5).Value2 = 9
$ws.Cells.Item($newRow, 6).Value2 = 100112044
$ws.Cells.Item($newRow, 7).Value2 = "Perejil"
$ws.Cells.Item($newRow, 8).Value2 = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value2 = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 30
$ws.Cells.Item($newRow, 11).Value2 = 4000
$ws.Cells.Item($newRow, 12).Value2 = 4000
$ws.Cells.Item($newRow, 13).Value2 = 4000
$ws.Cells.Item($newRow, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item($newRow, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item($newRow, 16).Value2 = 1333
$ws.Cells.Item($newRow, 17).Value2 = 3
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"
